# Update country data file: "Data" sheet -> "Summary" sheet, insert a
# "Source Type" title row, shift the MSME table down, and append a new
# "NSB" source-detail block at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename sheet -------------------------------------------------
$ws.Name = "Summary"

# --- 2. Start from a clean sheet and rebuild the layout --------------
# (the move shifts rows 5-9 down to 9-13 and adds new rows 7, 21, 22, so
# it's simplest/most reliable to clear and rewrite every cell fresh)
$ws.Cells.Clear()

# Row 1: country name ("name" style -> big 18pt font)
$c = $ws.Range("A1")
$c.Value = "Bhutan"
$c.Font.Size = 18

# Row 3: section title ("title" style -> bold)
$c = $ws.Range("A3")
$c.Value = "MSME Participation on the Economy"
$c.Font.Bold = $true

# Row 7: new sub-title ("title_" style -> bold + underline)
$c = $ws.Range("A7")
$c.Value = "Source Type: Statistical Institution (Most Widely Used)"
$c.Font.Bold = $true
$c.Font.Underline = $true

# Row 9: table header ("title" style -> bold)
$c = $ws.Range("B9")
$c.Value = "Micro"
$c.Font.Bold = $true

$c = $ws.Range("C9")
$c.Value = "SMEs"
$c.Font.Bold = $true

$c = $ws.Range("D9")
$c.Value = "MSMEs"
$c.Font.Bold = $true

# Row 10: Enterprises (absolute #)
# (data cells use a leading "'" so the numeric-looking text is kept as
# plain text -- matching the source workbook -- instead of being
# auto-converted to a number)
$c = $ws.Range("A10")
$c.Value = "Enterprises (absolute #)"
$c.Font.Bold = $true

$ws.Range("B10").Value = "'21210"
$ws.Range("C10").Value = "'3254"
$ws.Range("D10").Value = "'24464"

# Row 11: Enterprises density (per 1000 people)
$c = $ws.Range("A11")
$c.Value = "Enterprises density (per 1000 people)"
$c.Font.Bold = $true

$ws.Range("B11").Value = "'28.6"
$ws.Range("C11").Value = "'4.4"
$ws.Range("D11").Value = "'33"

# Row 12: Enterprises (% of total)
$c = $ws.Range("A12")
$c.Value = "Enterprises (% of total)"
$c.Font.Bold = $true

$ws.Range("B12").Value = "'86.2"
$ws.Range("C12").Value = "'13.2"
$ws.Range("D12").Value = "'99.4"

# Row 13: source footnote ("source" style -> italic)
$c = $ws.Range("A13")
$c.Value = "Source: NSB, 2012"
$c.Font.Italic = $true

# Row 21: source name ("title" style -> bold)
$c = $ws.Range("A21")
$c.Value = "NSB"
$c.Font.Bold = $true

# Row 22: full citation ("source" style -> italic)
$c = $ws.Range("A22")
$c.Value = "National Statistics Bureau (NSB), Royal Government of Bhutan, `"Statistical Yearbook of Bhutan 2013`", p. 118. Available at http://www.nsb.gov.bt/publication/files/pub9ot4338yv.pdf"
$c.Font.Italic = $true

Write-Host "Applied Bhutan Summary sheet update"
